$p = $ppt.ActivePresentation

# --- Slide 1: title text + font size change ---
$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1)
$title1.TextFrame.TextRange.Text = "Simple Two-Way Direct Loss Measurement Procedure"
$title1.TextFrame.TextRange.Font.Size = 36

# --- Slide 10: reposition/resize "Title 2" textbox and merge paragraphs ---
$s10 = $p.Slides.Item(10)
$title10 = $s10.Shapes.Item(1)
$title10.Left = 0
$title10.Top = 111502 / 12700
$title10.Width = 9144000 / 12700
$title10.Height = 731836 / 12700

# Merge the two paragraphs into a single line of text, dropping the explicit
# 32pt override so the run picks up the shape's default size.
$title10.TextFrame.TextRange.Text = "Link/P2P L2 Circuits - Counter-stamping in Hardware"
$title10.TextFrame.TextRange.Font.Size = 28

# --- Slide 5: give the footer placeholder an explicit position/size ---
$s5 = $p.Slides.Item(5)
$ftr5 = $s5.Shapes.Item(3)
$ftr5.Left = 3124200 / 12700
$ftr5.Top = 4786312 / 12700
$ftr5.Width = 2895600 / 12700
$ftr5.Height = 357188 / 12700
